$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 36: 602. Friend Requests II: Who Has the Most Friends
$ws.Range("A36").Value = "602. Friend Requests II: Who Has the Most Friends"
$ws.Range("B36").Value = "Medium"
$ws.Range("C36").Value = "Subqueries"
$ws.Range("D36").Value = "Use CTE as select id and count(*) as num from (select requester id UNION ALL accepter_id) group by id, then select id, num from cte where num = (select max(num) from cte)"
$ws.Range("E36").Value = "https://leetcode.com/problems/friend-requests-ii-who-has-the-most-friends/solutions/3550546/simple-solution/?envType=study-plan-v2&envId=top-sql-50 "

# Row 37: 585. Investments in 2016 (note: E filled before D to match original string pool order)
$ws.Range("A37").Value = "585. Investments in 2016"
$ws.Range("B37").Value = "Medium"
$ws.Range("C37").Value = "Subqueries"
$ws.Range("E37").Value = "https://leetcode.com/problems/investments-in-2016/solutions/4156632/880ms-runtime-beat-97-users-0mb-memory/?envType=study-plan-v2&envId=top-sql-50 "
$ws.Range("D37").Value = "We can either use subqueries, or CTE with window functions."

# Match the styling of the existing Medium rows (B column fill)
$ws.Range("B36").Interior.Color = $ws.Range("B35").Interior.Color
$ws.Range("B37").Interior.Color = $ws.Range("B35").Interior.Color

# Add hyperlinks for E36 and E37
$ws.Hyperlinks.Add($ws.Range("E36"), "https://leetcode.com/problems/friend-requests-ii-who-has-the-most-friends/solutions/3550546/simple-solution/?envType=study-plan-v2&envId=top-sql-50 ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E37"), "https://leetcode.com/problems/investments-in-2016/solutions/4156632/880ms-runtime-beat-97-users-0mb-memory/?envType=study-plan-v2&envId=top-sql-50 ") | Out-Null

# Grow the table (Table2) to include the two new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E37"))

# Update selection to D41, per diff
$ws.Range("D41").Select()
